# Update the "想去人数" (F column) figures across the four sheets to match
# the freshly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 7603   # 上海·艺术与潮流·遇见EVA 中国首展
$ws1.Range("F11").Value = 6      # 上海·X-party国漫游戏嘉年华02（免费活动）
$ws1.Range("F12").Value = 20     # 上海·次元码头 舞蹈嘉年华~宅舞精英赛（免费活动）
$ws1.Range("F19").Value = 1394   # 上海 洛天依歌行宇宙·无限遨游 沉浸式体验展
$ws1.Range("F24").Value = 4157   # 上海·原神ONLY逐月节·原神&崩铁&崩三&绝区零·同人动漫嘉年华
$ws1.Range("F25").Value = 3316   # 上海·城市动漫节2th
$ws1.Range("F26").Value = 272    # 上海·夜蓝诗·恋与深空同人only
$ws1.Range("F27").Value = 83     # 上海·宫村优子粉丝见面会
$ws1.Range("F28").Value = 83     # 上海·宫村优子粉丝见面会
$ws1.Range("F48").Value = 1970   # 上海·次元同人动漫节only·旅行盛宴3.0

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 112    # 上海·日本次世代神秘创作歌手 和ぬか(wanuka) 演出
$ws2.Range("F10").Value = 594    # 上海·神山羊2024巡演ENCOUNTER
$ws2.Range("F27").Value = 4449   # 上海·洛天依2024无限共鸣演唱会
$ws2.Range("F28").Value = 4449   # 上海·洛天依2024无限共鸣演唱会

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value  = 3013   # 上海·「BanG Dream! It's MyGO!!!!! × animate cafe」
$ws3.Range("F13").Value = 2025   # 上海·东方明珠·「光与夜之恋...」线条大作战主题店
$ws3.Range("F14").Value = 8700   # 上海·大悦城·「光与夜之恋...」线条大作战主题餐厅
$ws3.Range("F15").Value = 851    # 上海·「HUNTER×HUNTER × animate cafe」

# --- Sheet 4: 全部类型 (All types, aggregate of the above) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 3013   # 上海·「BanG Dream! It's MyGO!!!!! × animate cafe」
$ws4.Range("F9").Value  = 7603   # 上海·艺术与潮流·遇见EVA 中国首展
$ws4.Range("F14").Value = 112    # 上海·日本次世代神秘创作歌手 和ぬか(wanuka) 演出
$ws4.Range("F15").Value = 20     # 上海·次元码头 舞蹈嘉年华~宅舞精英赛（免费活动）
$ws4.Range("F17").Value = 851    # 上海·「HUNTER×HUNTER × animate cafe」
$ws4.Range("F19").Value = 594    # 上海·神山羊2024巡演ENCOUNTER
$ws4.Range("F20").Value = 594    # 上海·神山羊2024巡演ENCOUNTER
$ws4.Range("F27").Value = 1394   # 上海 洛天依歌行宇宙·无限遨游 沉浸式体验展
$ws4.Range("F31").Value = 4157   # 上海·原神ONLY逐月节·原神&崩铁&崩三&绝区零·同人动漫嘉年华
$ws4.Range("F32").Value = 3316   # 上海·城市动漫节2th
$ws4.Range("F33").Value = 272    # 上海·夜蓝诗·恋与深空同人only
$ws4.Range("F34").Value = 83     # 上海·宫村优子粉丝见面会
$ws4.Range("F35").Value = 83     # 上海·宫村优子粉丝见面会
$ws4.Range("F48").Value = 1970   # 上海·次元同人动漫节only·旅行盛宴3.0
$ws4.Range("F50").Value = 4449   # 上海·洛天依2024无限共鸣演唱会
